$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CellAddress (H), Attribute (E), TagType (I) and ProcessFilter (J)
# columns for rows 2-7 to their new values.

$ws.Range("E2").Value = "ACT_BND"
$ws.Range("H2").Value = "3,1"
$ws.Range("I2").Value = "TFM_INS"
$ws.Range("J2").Value = '{"pset_pn": "IMP*GZ", "pset_set": "IRE"}'

$ws.Range("E3").Value = "ACT_BND"
$ws.Range("H3").Value = "4,1"
$ws.Range("I3").Value = "TFM_INS"
$ws.Range("J3").Value = '{"pset_pn": "IMPDEMZ", "pset_set": "IRE"}'

$ws.Range("E4").Value = "ACT_BND"
$ws.Range("H4").Value = "5,1"
$ws.Range("I4").Value = "TFM_INS"
$ws.Range("J4").Value = '{"pset_pn": "IMP*Z", "pset_set": "IRE"}'

$ws.Range("E5").Value = "PRC_TSL"
$ws.Range("H5").Value = "9,1"
$ws.Range("I5").Value = "TFM_INS-TXT"
$ws.Range("J5").Value = '{"pset_pn": "IMP*Z", "pset_set": "IRE"}'

$ws.Range("E6").Value = "ACTCOST"
$ws.Range("H6").Value = "22,1"
$ws.Range("I6").Value = "TFM_UPD"
$ws.Range("J6").Value = '{"pset_pn": "IMP*Z", "pset_set": "IRE"}'

$ws.Range("E7").Value = "ACTCOST"
$ws.Range("H7").Value = "23,1"
$ws.Range("I7").Value = "TFM_UPD"
$ws.Range("J7").Value = '{"pset_pn": "IMPDEMZ", "pset_set": "IRE"}'
